$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: rename existing strings in place (same cells, new text)
$ws.Range("A1").Value = "sulakshana"
$ws.Range("B1").Value = "sulaa"
$ws.Range("C1").Value = "ncc"
$ws.Range("D1").Value = 23
$ws.Range("E1").Value = "m"

# Row 2: new strings + updated number
$ws.Range("A2").Value = "niro"
$ws.Range("B2").Value = "ssssss"
$ws.Range("C2").Value = "dsc"
$ws.Range("D2").Value = 24
$ws.Range("E2").Value = "m"
